$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

Set-TextValue "D2" "29.005.99"
Set-TextValue "E2" "  -0.32%  "

Set-TextValue "D3" "1.918.57"
Set-TextValue "E3" "  +0.33%  "

Set-TextValue "E4" "  -0.44%  "

Set-TextValue "D5" "325.18"
Set-TextValue "E5" "  +0.06%  "

Set-TextValue "E6" "  -0.42%  "

Set-TextValue "D7" "0.4601"
Set-TextValue "E7" "  -0.05%  "

Set-TextValue "D8" "0.3825"
Set-TextValue "E8" "  -0.19%  "

Set-TextValue "D9" "0.07706"
Set-TextValue "E9" "  -0.70%  "

Set-TextValue "D10" "0.9816"
Set-TextValue "E10" "  +1.72%  "

Set-TextValue "D11" "22.22"
Set-TextValue "E11" "  +0.48%  "

Set-TextValue "D12" "1.894.85"
Set-TextValue "E12" "  -1.02%  "

Set-TextValue "D13" "5.693"
Set-TextValue "E13" "  -0.23%  "

Set-TextValue "D14" "6.970"
Set-TextValue "E14" "  -0.63%  "

Set-TextValue "D15" "0.06995"
Set-TextValue "E15" "  -1.24%  "

Set-TextValue "E16" "  -0.44%  "

Set-TextValue "D17" "84.24"
Set-TextValue "E17" "  -0.52%  "

Set-TextValue "D18" "0.000009490"
Set-TextValue "E18" "  -1.12%  "

Set-TextValue "E19" "  -1.40%  "

Set-TextValue "E20" "  -0.39%  "

Set-TextValue "D21" "28.995.15"
Set-TextValue "E21" "  -0.32%  "

Set-TextValue "D22" "5.340"
Set-TextValue "E22" "  -2.01%  "

Set-TextValue "D23" "10.96"
Set-TextValue "E23" "  +0.08%  "

Set-TextValue "D24" "2.157.41"
Set-TextValue "E24" "  +0.02%  "

Set-TextValue "D25" "2.090"
Set-TextValue "E25" "  -0.33%  "

Set-TextValue "D26" "158.47"
Set-TextValue "E26" "  +0.68%  "

Set-TextValue "E27" "  -0.85%  "

Set-TextValue "D28" "5.705"
Set-TextValue "E28" "  +0.57%  "

Set-TextValue "D29" "117.82"
Set-TextValue "E29" "  +0.17%  "

Set-TextValue "D30" "1.864"
Set-TextValue "E30" "  +2.41%  "

Set-TextValue "D31" "0.09316"
Set-TextValue "E31" "  +0.09%  "

Set-TextValue "D32" "0.8684"
Set-TextValue "E32" "  +1.63%  "

Set-TextValue "D33" "5.116"
Set-TextValue "E33" "  +0.35%  "

Set-TextValue "D34" "1.253"
Set-TextValue "E34" "  -0.32%  "

Set-TextValue "D35" "3.048"
Set-TextValue "E35" "  -1.16%  "

Set-TextValue "D36" "0.05711"
Set-TextValue "E36" "  +0.30%  "

Set-TextValue "D37" "1.156"
Set-TextValue "E37" "  -0.21%  "

Set-TextValue "D38" "1.000"
Set-TextValue "E38" "  -0.52%  "

Set-TextValue "E39" "  -0.52%  "

Set-TextValue "D40" "3.050"
Set-TextValue "E40" "  +12.34%  "

Set-TextValue "D41" "7.541"
Set-TextValue "E41" "  +0.26%  "

Set-TextValue "D42" "0.5514"
Set-TextValue "E42" "  -0.93%  "

Set-TextValue "D43" "0.1751"
Set-TextValue "E43" "  -0.52%  "

Set-TextValue "D44" "9.394"
Set-TextValue "E44" "  +2.00%  "

Set-TextValue "D45" "0.000002882"
Set-TextValue "E45" "  -0.07%  "

Set-TextValue "D46" "2.193"
Set-TextValue "E46" "  +6.45%  "

Set-TextValue "D47" "0.5194"
Set-TextValue "E47" "  -0.55%  "

Set-TextValue "D48" "11.20"
Set-TextValue "E48" "  -0.60%  "

Set-TextValue "D49" "0.06909"
Set-TextValue "E49" "  +1.33%  "

Set-TextValue "E50" "  -0.45%  "

Set-TextValue "D51" "110.42"
Set-TextValue "E51" "  -0.32%  "
